$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap "Polonia" and "Banglades" labels (rows 33 and 34) ---
$ws.Range("A33").Value = "Banglades"
$ws.Range("A34").Value = "Polonia"

# --- Row 33 (Banglades, updated stats) ---
$ws.Range("B33").Value = 17822
$ws.Range("C33").Value = 1162
$ws.Range("D33").Value = 3361
$ws.Range("E33").Value = 14192
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 19
$ws.Range("H33").Value = 269

# --- Row 34 (Polonia, updated stats) ---
$ws.Range("B34").Value = 17062
$ws.Range("C34").Value = 141
$ws.Range("D34").Value = 6410
$ws.Range("E34").Value = 9805
$ws.Range("F34").Value = 160
$ws.Range("G34").Value = 8
$ws.Range("H34").Value = 847

# --- Row 35 (Israel, updated stats) ---
$ws.Range("B35").Value = 16539
$ws.Range("C35").Value = 10
$ws.Range("D35").Value = 12173
$ws.Range("E35").Value = 4104
$ws.Range("F35").Value = 61
$ws.Range("G35").Value = 2
$ws.Range("H35").Value = 262

# --- Row 40 (Indonesia, updated stats) ---
$ws.Range("B40").Value = 15438
$ws.Range("C40").Value = 689
$ws.Range("D40").Value = 3287
$ws.Range("E40").Value = 11123
$ws.Range("G40").Value = 21
$ws.Range("H40").Value = 1028

# --- Row 54 (Malasia, updated stats) ---
$ws.Range("B54").Value = 6779
$ws.Range("C54").Value = 37
$ws.Range("D54").Value = 5281
$ws.Range("E54").Value = 1387
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 111

# --- Row 97 (Hong Kong, updated stats) ---
$ws.Range("B97").Value = 1051
$ws.Range("C97").Value = 3
$ws.Range("D97").Value = 1008
$ws.Range("E97").Value = 39

# --- Row 115 (Tayikistan, updated stats) ---
$ws.Range("E115").Value = 706
$ws.Range("G115").Value = 2
$ws.Range("H115").Value = 23

# --- Row 136 (Benin, updated stats) ---
$ws.Range("D136").Value = 83
$ws.Range("E136").Value = 242

# --- Update "Datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 11:05"
